$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2" = 37.98277566666666
    "H2" = 113.948327
    "I2" = 0.697850645410475
    "J2" = 0.6978506454104751
    "M2" = 449.3583473333333
    "N2" = 1348.075042
    "O2" = 0.959704436884883
    "P2" = 0.9597044368848828
    "Q2" = 17067.87730070608
    "R2" = 153610.8957063547
    "S2" = 0.6697303606834122
    "T2" = 0.6697303606834121
    "G3" = 37.98277566666666
    "H3" = 113.948327
    "I3" = 0.697850645410475
    "J3" = 0.6978506454104751
    "O3" = 0.01202662913387072
    "P3" = 0.01202662913387072
    "Q3" = 213.8877580521434
    "R3" = 1924.989822469291
    "S3" = 0.008392790903184103
    "T3" = 0.008392790903184103
    "G4" = 37.98277566666666
    "H4" = 113.948327
    "I4" = 0.697850645410475
    "J4" = 0.6978506454104751
    "M4" = 7.708291333333332
    "N4" = 23.124874
    "O4" = 0.01646276615823874
    "P4" = 0.01646276615823874
    "Q4" = 292.7823004873108
    "R4" = 2635.040704385798
    "S4" = 0.01148855198876863
    "T4" = 0.01148855198876863
    "G5" = 37.98277566666666
    "H5" = 113.948327
    "I5" = 0.697850645410475
    "J5" = 0.6978506454104751
    "M5" = 1.356257333333333
    "N5" = 4.068772
    "O5" = 0.002896588408965574
    "P5" = 0.002896588408965573
    "Q5" = 51.51441803827155
    "R5" = 463.629762344444
    "S5" = 0.002021386090685127
    "T5" = 0.002021386090685127
    "G6" = 37.98277566666666
    "H6" = 113.948327
    "I6" = 0.697850645410475
    "J6" = 0.6978506454104751
    "M6" = 4.171694666666667
    "N6" = 12.515084
    "O6" = 0.008909579414042005
    "P6" = 0.008909579414042003
    "Q6" = 158.4525426738298
    "R6" = 1426.072884064468
    "S6" = 0.006217555744425095
    "T6" = 0.006217555744425095
    "I7" = 0.1779541659542351
    "J7" = 0.1779541659542352
    "M7" = 449.3583473333333
    "N7" = 1348.075042
    "O7" = 0.959704436884883
    "P7" = 0.9597044368848828
    "Q7" = 4352.363775303003
    "R7" = 39171.27397772703
    "S7" = 0.1707834026284283
    "T7" = 0.1707834026284283
    "I8" = 0.1779541659542351
    "J8" = 0.1779541659542352
    "O8" = 0.01202662913387072
    "P8" = 0.01202662913387072
    "S8" = 0.002140188756758869
    "T8" = 0.002140188756758869
    "I9" = 0.1779541659542351
    "J9" = 0.1779541659542352
    "M9" = 7.708291333333332
    "N9" = 23.124874
    "O9" = 0.01646276615823874
    "P9" = 0.01646276615823874
    "Q9" = 74.6604311854371
    "R9" = 671.9438806689341
    "S9" = 0.002929617820988983
    "T9" = 0.002929617820988983
    "I10" = 0.1779541659542351
    "J10" = 0.1779541659542352
    "M10" = 1.356257333333333
    "N10" = 4.068772
    "O10" = 0.002896588408965574
    "P10" = 0.002896588408965573
    "Q10" = 13.13634279327245
    "R10" = 118.227085139452
    "S10" = 0.0005154599744301737
    "T10" = 0.0005154599744301736
    "I11" = 0.1779541659542351
    "J11" = 0.1779541659542352
    "M11" = 4.171694666666667
    "N11" = 12.515084
    "O11" = 0.008909579414042005
    "P11" = 0.008909579414042003
    "Q11" = 40.40590957433823
    "R11" = 363.6531861690441
    "S11" = 0.001585496773628868
    "T11" = 0.001585496773628868
    "G12" = 0.5676613333333332
    "H12" = 1.702984
    "I12" = 0.01042953867610283
    "J12" = 0.01042953867610283
    "M12" = 449.3583473333333
    "N12" = 1348.075042
    "O12" = 0.959704436884883
    "P12" = 0.9597044368848828
    "Q12" = 255.083358591703
    "R12" = 2295.750227325328
    "S12" = 0.01000927454211837
    "T12" = 0.01000927454211837
    "G13" = 0.5676613333333332
    "H13" = 1.702984
    "I13" = 0.01042953867610283
    "J13" = 0.01042953867610283
    "O13" = 0.01202662913387072
    "P13" = 0.01202662913387072
    "Q13" = 3.196601822496889
    "R13" = 28.769416402472
    "S13" = 0.0001254321936948497
    "T13" = 0.0001254321936948497
    "G14" = 0.5676613333333332
    "H14" = 1.702984
    "I14" = 0.01042953867610283
    "J14" = 0.01042953867610283
    "M14" = 7.708291333333332
    "N14" = 23.124874
    "O14" = 0.01646276615823874
    "P14" = 0.01646276615823874
    "Q14" = 4.375698936001776
    "R14" = 39.38129042401599
    "S14" = 0.0001716990563629878
    "T14" = 0.0001716990563629877
    "G15" = 0.5676613333333332
    "H15" = 1.702984
    "I15" = 0.01042953867610283
    "J15" = 0.01042953867610283
    "M15" = 1.356257333333333
    "N15" = 4.068772
    "O15" = 0.002896588408965574
    "P15" = 0.002896588408965573
    "Q15" = 0.769894846183111
    "R15" = 6.929053615648
    "S15" = 0.00003021008084005762
    "T15" = 0.00003021008084005762
    "G16" = 0.5676613333333332
    "H16" = 1.702984
    "I16" = 0.01042953867610283
    "J16" = 0.01042953867610283
    "M16" = 4.171694666666667
    "N16" = 12.515084
    "O16" = 0.008909579414042005
    "P16" = 0.008909579414042003
    "Q16" = 2.368109756739555
    "R16" = 21.312987810656
    "S16" = 0.00009292280308656067
    "T16" = 0.00009292280308656067
    "G17" = 5.823095333333334
    "H17" = 17.469286
    "I17" = 0.1069866739681064
    "J17" = 0.1069866739681064
    "M17" = 449.3583473333333
    "N17" = 1348.075042
    "O17" = 0.959704436884883
    "P17" = 0.9597044368848828
    "Q17" = 2616.656495351112
    "R17" = 23549.90845816001
    "S17" = 0.1026755856947482
    "T17" = 0.1026755856947481
    "G18" = 5.823095333333334
    "H18" = 17.469286
    "I18" = 0.1069866739681064
    "J18" = 0.1069866739681064
    "O18" = 0.01202662913387072
    "P18" = 0.01202662913387072
    "Q18" = 32.79088439193756
    "R18" = 295.117959527438
    "S18" = 0.001286689050080757
    "T18" = 0.001286689050080757
    "G19" = 5.823095333333334
    "H19" = 17.469286
    "I19" = 0.1069866739681064
    "J19" = 0.1069866739681064
    "M19" = 7.708291333333332
    "N19" = 23.124874
    "O19" = 0.01646276615823874
    "P19" = 0.01646276615823874
    "Q19" = 44.88611529110711
    "R19" = 403.975037619964
    "S19" = 0.001761296595584664
    "T19" = 0.001761296595584664
    "G20" = 5.823095333333334
    "H20" = 17.469286
    "I20" = 0.1069866739681064
    "J20" = 0.1069866739681064
    "M20" = 1.356257333333333
    "N20" = 4.068772
    "O20" = 0.002896588408965574
    "P20" = 0.002896588408965573
    "Q20" = 7.897615748532446
    "R20" = 71.078541736792
    "S20" = 0.0003098963597297961
    "T20" = 0.000309896359729796
    "G21" = 5.823095333333334
    "H21" = 17.469286
    "I21" = 0.1069866739681064
    "J21" = 0.1069866739681064
    "M21" = 4.171694666666667
    "N21" = 12.515084
    "O21" = 0.008909579414042005
    "P21" = 0.008909579414042003
    "Q21" = 24.29217574555823
    "R21" = 218.629581710024
    "S21" = 0.0009532062679630647
    "T21" = 0.0009532062679630646
    "G22" = 0.3689676666666666
    "H22" = 1.106903
    "I22" = 0.006778975991080511
    "J22" = 0.006778975991080512
    "M22" = 449.3583473333333
    "N22" = 1348.075042
    "O22" = 0.959704436884883
    "P22" = 0.9597044368848828
    "Q22" = 165.7987009127695
    "R22" = 1492.188308214926
    "S22" = 0.006505813336176064
    "T22" = 0.006505813336176063
    "G23" = 0.3689676666666666
    "H23" = 1.106903
    "I23" = 0.006778975991080511
    "J23" = 0.006778975991080512
    "O23" = 0.01202662913387072
    "P23" = 0.01202662913387072
    "Q23" = 2.077722484255444
    "R23" = 18.699502358299
    "S23" = 0.00008152823015213899
    "T23" = 0.00008152823015213899
    "G24" = 0.3689676666666666
    "H24" = 1.106903
    "I24" = 0.006778975991080511
    "J24" = 0.006778975991080512
    "M24" = 7.708291333333332
    "N24" = 23.124874
    "O24" = 0.01646276615823874
    "P24" = 0.01646276615823874
    "Q24" = 2.844110267246888
    "R24" = 25.596992405222
    "S24" = 0.0001116006965334732
    "T24" = 0.0001116006965334732
    "G25" = 0.3689676666666666
    "H25" = 1.106903
    "I25" = 0.006778975991080511
    "J25" = 0.006778975991080512
    "M25" = 1.356257333333333
    "N25" = 4.068772
    "O25" = 0.002896588408965574
    "P25" = 0.002896588408965573
    "Q25" = 0.5004151036795556
    "R25" = 4.503735933116
    "S25" = 0.00001963590328041972
    "T25" = 0.00001963590328041972
    "G26" = 0.3689676666666666
    "H26" = 1.106903
    "I26" = 0.006778975991080511
    "J26" = 0.006778975991080512
    "M26" = 4.171694666666667
    "N26" = 12.515084
    "O26" = 0.008909579414042005
    "P26" = 0.008909579414042003
    "Q26" = 1.539220447205778
    "R26" = 13.852984024852
    "S26" = 0.00006039782493841592
    "T26" = 0.00006039782493841591
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
